$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.694.48"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "'1.902.88"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'312.29"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.5211"
$ws.Range("E7").Value = "  +7.55%  "
$ws.Range("D8").Value = "'0.3773"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "'0.07234"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("E10").Value = "  +3.43%  "
$ws.Range("D11").Value = "'0.8949"
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Value = "'0.07619"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "'1.891.53"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "'5.437"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "'91.95"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'0.000008704"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'27.733.14"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'14.43"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "'2.124.18"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'6.573"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "'153.41"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'1.873"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.28"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.150"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").Value = "'114.44"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").Value = "'4.838"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").Value = "'0.09009"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'4.880"
$ws.Range("E32").Value = "  +5.30%  "
$ws.Range("D33").Value = "'3.175"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("D34").Value = "'1.231"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").Value = "'0.7671"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "'2.622"
$ws.Range("E36").Value = "  +3.56%  "
$ws.Range("D37").Value = "'0.02080"
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("D38").Value = "'3.063"
$ws.Range("E38").Value = "  +2.78%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.05276"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5481"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "'6.637"
$ws.Range("E42").Value = "  -4.09%  "
$ws.Range("D43").Value = "'114.12"
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("D44").Value = "'8.488"
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").Value = "'0.1506"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").Value = "'0.4778"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "'10.43"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "'0.9995"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'1.612"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").Value = "'66.40"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "'0.05990"
$ws.Range("E51").Value = "  -1.06%  "
